$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 273.33334
$ws.Range("I6").Value = 305
$ws.Range("K6").Value = 915
$ws.Range("M6").Value = -803
$ws.Range("H13").Value = 17999.5
$ws.Range("J13").Value = 17999.5
$ws.Range("L13").Value = 17999.5
$ws.Range("N13").Value = -18337.5
$ws.Range("H15").Value = 1204764.2
$ws.Range("I15").Value = 1204764.2
$ws.Range("K15").Value = 3614292.6
$ws.Range("M15").Value = -3614123.6
$ws.Range("H39").Value = 1979
$ws.Range("I39").Value = 808.63635
$ws.Range("J39").Value = 3051.8333
$ws.Range("K39").Value = 2425.90905
$ws.Range("L39").Value = 9155.499899999999
$ws.Range("M39").Value = -2129.90905
$ws.Range("N39").Value = -9747.499899999999
$ws.Range("H62").Value = 3999.5
$ws.Range("I62").Value = 3999
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 3999
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -3375
$ws.Range("N62").Value = -5248
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").ClearContents()
$ws.Range("H65").Value = 3999.5
$ws.Range("I65").Value = 3999
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 19995
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -16875
$ws.Range("N65").Value = -26240
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").ClearContents()
$ws.Range("H107").Value = 1122.25
$ws.Range("I107").Value = 1122.25
$ws.Range("K107").Value = 1122.25
$ws.Range("M107").Value = 797.75
$ws.Range("H116").Value = 5161.722
$ws.Range("I116").Value = 5000.8
$ws.Range("J116").Value = 5966.3335
$ws.Range("K116").Value = 5000.8
$ws.Range("L116").Value = 5966.3335
$ws.Range("M116").Value = -1558.8
$ws.Range("N116").Value = -12850.3335
$ws.Range("H125").Value = 11838762
$ws.Range("J125").Value = 13901167
$ws.Range("L125").Value = 125110503
$ws.Range("N125").Value = -125115423
$ws.Range("H131").Value = 874.75
$ws.Range("I131").Value = 874.75
$ws.Range("K131").Value = 2624.25
$ws.Range("M131").Value = 2415.75
$ws.Range("H132").Value = 3685.2354
$ws.Range("I132").Value = 3685.2354
$ws.Range("K132").Value = 11055.7062
$ws.Range("M132").Value = -8525.706200000001
$ws.Range("H133").Value = 95000
$ws.Range("J133").Value = 95000
$ws.Range("L133").Value = 95000
$ws.Range("N133").Value = -105120
$ws.Range("H135").Value = 187500820
$ws.Range("I135").Value = 83333930
$ws.Range("K135").Value = 750005370
$ws.Range("M135").Value = -750002835
$ws.Range("H137").Value = 2866.5454
$ws.Range("I137").Value = 2564.6667
$ws.Range("J137").Value = 4225
$ws.Range("K137").Value = 7694.000100000001
$ws.Range("L137").Value = 12675
$ws.Range("M137").Value = -5144.000100000001
$ws.Range("N137").Value = -17775
$ws.Range("H138").Value = 2594.6316
$ws.Range("I138").Value = 1520.8462
$ws.Range("J138").Value = 3153
$ws.Range("K138").Value = 4562.5386
$ws.Range("L138").Value = 9459
$ws.Range("M138").Value = 577.4614000000001
$ws.Range("N138").Value = -19739
$ws.Range("H141").Value = 6006.6
$ws.Range("I141").Value = 5924.5557
$ws.Range("K141").Value = 17773.6671
$ws.Range("M141").Value = -12593.6671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3723.8125
$ws.Range("I32").Value = 1946.3928
$ws.Range("J32").Value = 16165.75
$ws.Range("K32").Value = 1946.3928
$ws.Range("L32").Value = 16165.75
$ws.Range("M32").Value = -1659.3928
$ws.Range("N32").Value = -16739.75
$ws.Range("H45").Value = 5783.9287
$ws.Range("I45").Value = 7077.9
$ws.Range("J45").Value = 2549
$ws.Range("K45").Value = 7077.9
$ws.Range("L45").Value = 2549
$ws.Range("M45").Value = -6700.9
$ws.Range("N45").Value = -3303
$ws.Range("H97").Value = 660.8889
$ws.Range("I97").Value = 700.5833
$ws.Range("K97").Value = 700.5833
$ws.Range("M97").Value = -204.5833
$ws.Range("H119").Value = 15000
$ws.Range("J119").Value = 15000
$ws.Range("L119").Value = 15000
$ws.Range("N119").Value = -24676
$ws.Range("H122").Value = 3386.9707
$ws.Range("I122").Value = 2539.9312
$ws.Range("K122").Value = 7619.7936
$ws.Range("M122").Value = -5169.7936
$ws.Range("H132").Value = 4169324.8
$ws.Range("I132").Value = 4547900
$ws.Range("K132").Value = 13643700
$ws.Range("M132").Value = -13641170

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2270.682
$ws.Range("J20").Value = 1963.5555
$ws.Range("L20").Value = 1963.5555
$ws.Range("N20").Value = -2457.5555
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H86").Value = 1634.4889
$ws.Range("I86").Value = 1620.3429
$ws.Range("J86").Value = 1684
$ws.Range("K86").Value = 1620.3429
$ws.Range("L86").Value = 1684
$ws.Range("M86").Value = -497.3429000000001
$ws.Range("N86").Value = -3930
$ws.Range("H89").Value = 1634.4889
$ws.Range("I89").Value = 1620.3429
$ws.Range("J89").Value = 1684
$ws.Range("K89").Value = 8101.7145
$ws.Range("L89").Value = 8420
$ws.Range("M89").Value = -2485.7145
$ws.Range("N89").Value = -19652
$ws.Range("H105").Value = 2635.9443
$ws.Range("I105").Value = 2111.3076
$ws.Range("K105").Value = 2111.3076
$ws.Range("M105").Value = -364.3076000000001
$ws.Range("H107").Value = 147201.58
$ws.Range("I107").Value = 5068.6665
$ws.Range("K107").Value = 5068.6665
$ws.Range("M107").Value = -3148.6665
$ws.Range("H134").Value = 50006096
$ws.Range("I134").Value = 50006096
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 150018288
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -150015753
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7336.4287
$ws.Range("I31").Value = 4968.522
$ws.Range("K31").Value = 4968.522
$ws.Range("M31").Value = -4673.522
$ws.Range("H32").Value = 4533
$ws.Range("I32").Value = 4533
$ws.Range("K32").Value = 4533
$ws.Range("M32").Value = -4217
$ws.Range("H34").Value = 7336.4287
$ws.Range("I34").Value = 4968.522
$ws.Range("K34").Value = 4968.522
$ws.Range("M34").Value = -4766.522
$ws.Range("H45").Value = 1000
$ws.Range("I45").Value = 1000
$ws.Range("K45").Value = 1000
$ws.Range("M45").Value = -407
$ws.Range("H86").Value = 14128.286
$ws.Range("J86").Value = 15469.333
$ws.Range("L86").Value = 15469.333
$ws.Range("N86").Value = -17715.333
$ws.Range("H88").Value = 17499.875
$ws.Range("J88").Value = 19999.75
$ws.Range("L88").Value = 19999.75
$ws.Range("N88").Value = -20811.75
$ws.Range("H89").Value = 14128.286
$ws.Range("J89").Value = 15469.333
$ws.Range("L89").Value = 77346.66500000001
$ws.Range("N89").Value = -88578.66500000001
$ws.Range("H91").Value = 17499.875
$ws.Range("J91").Value = 19999.75
$ws.Range("L91").Value = 19999.75
$ws.Range("N91").Value = -22807.75
$ws.Range("H99").Value = 3050.3125
$ws.Range("I99").Value = 2953.9285
$ws.Range("K99").Value = 2953.9285
$ws.Range("M99").Value = -1455.9285
$ws.Range("H126").Value = 3050.3125
$ws.Range("I126").Value = 2953.9285
$ws.Range("K126").Value = 8861.7855
$ws.Range("M126").Value = -6391.7855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5829.091
$ws.Range("I3").Value = 4912.1
$ws.Range("K3").Value = 14736.3
$ws.Range("M3").Value = -14624.3
$ws.Range("H4").Value = 193029.44
$ws.Range("I4").Value = 213160.23
$ws.Range("K4").Value = 639480.6900000001
$ws.Range("M4").Value = -639368.6900000001
$ws.Range("H92").Value = 931.3333
$ws.Range("I92").Value = 947.5
$ws.Range("J92").Value = 899
$ws.Range("K92").Value = 2842.5
$ws.Range("L92").Value = 2697
$ws.Range("M92").Value = -1594.5
$ws.Range("N92").Value = -5193
$ws.Range("H99").Value = 5912.25
$ws.Range("H109").Value = 1426
$ws.Range("I109").Value = 1229.375
$ws.Range("J109").Value = 2999
$ws.Range("K109").Value = 3688.125
$ws.Range("L109").Value = 8997
$ws.Range("M109").Value = -2648.125
$ws.Range("N109").Value = -11077
$ws.Range("H110").Value = 20499.5
$ws.Range("I110").Value = 20000
$ws.Range("J110").Value = 20999
$ws.Range("K110").Value = 60000
$ws.Range("L110").Value = 62997
$ws.Range("M110").Value = -55910
$ws.Range("N110").Value = -71177
$ws.Range("H113").Value = 67748.60000000001
$ws.Range("J113").Value = 1253.2
$ws.Range("L113").Value = 3759.6
$ws.Range("N113").Value = -8099.6
$ws.Range("H114").Value = 144221.28
$ws.Range("J114").Value = 1810
$ws.Range("L114").Value = 5430
$ws.Range("N114").Value = -11938
$ws.Range("H117").Value = 2767.7
$ws.Range("J117").Value = 3359.75
$ws.Range("L117").Value = 10079.25
$ws.Range("N117").Value = -16963.25
$ws.Range("H120").Value = 575
$ws.Range("I120").Value = 575
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 1725
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = 3113
$ws.Range("N120").ClearContents()
$ws.Range("H129").Value = 2126.7222
$ws.Range("I129").Value = 1893.6666
$ws.Range("J129").Value = 2359.7778
$ws.Range("K129").Value = 5680.9998
$ws.Range("L129").Value = 7079.3334
$ws.Range("M129").Value = -680.9997999999996
$ws.Range("N129").Value = -17079.3334
$ws.Range("H131").Value = 1612.7727
$ws.Range("I131").Value = 1165.4166
$ws.Range("J131").Value = 2149.6
$ws.Range("K131").Value = 3496.2498
$ws.Range("L131").Value = 6448.799999999999
$ws.Range("M131").Value = 1543.7502
$ws.Range("N131").Value = -16528.8
$ws.Range("H139").Value = 1433.2693
$ws.Range("I139").Value = 1436.0416
$ws.Range("K139").Value = 4308.1248
$ws.Range("M139").Value = 831.8752000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 4325
$ws.Range("I31").Value = 4325
$ws.Range("K31").Value = 4325
$ws.Range("M31").Value = -4033
$ws.Range("H37").Value = 4325
$ws.Range("I37").Value = 4325
$ws.Range("K37").Value = 4325
$ws.Range("M37").Value = -4048
$ws.Range("H80").Value = 1999
$ws.Range("I80").Value = 1999
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 1999
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -1001
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 1999
$ws.Range("I83").Value = 1999
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 9995
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -5003
$ws.Range("N83").ClearContents()
$ws.Range("H102").Value = 3000.4348
$ws.Range("I102").Value = 3000.4348
$ws.Range("K102").Value = 3000.4348
$ws.Range("M102").Value = -1378.4348
$ws.Range("H107").Value = 1737.3334
$ws.Range("I107").Value = 1552.3529
$ws.Range("J107").Value = 2523.5
$ws.Range("K107").Value = 1552.3529
$ws.Range("L107").Value = 2523.5
$ws.Range("M107").Value = 367.6470999999999
$ws.Range("N107").Value = -6363.5
$ws.Range("H122").Value = 116727.82
$ws.Range("I122").Value = 174999.72
$ws.Range("J122").Value = 14752
$ws.Range("K122").Value = 524999.16
$ws.Range("L122").Value = 44256
$ws.Range("M122").Value = -522549.16
$ws.Range("N122").Value = -49156
$ws.Range("H126").Value = 3035.7778
$ws.Range("I126").Value = 2927.75
$ws.Range("K126").Value = 8783.25
$ws.Range("M126").Value = -6313.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4051.5715
$ws.Range("I7").Value = 4051.5715
$ws.Range("K7").Value = 4051.5715
$ws.Range("M7").Value = -3939.5715
$ws.Range("H16").Value = 735.2857
$ws.Range("I16").Value = 735.2857
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 735.2857
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -565.2857
$ws.Range("N16").ClearContents()
$ws.Range("H40").Value = 2199.353
$ws.Range("I40").Value = 1901.091
$ws.Range("J40").Value = 2746.1667
$ws.Range("K40").Value = 1901.091
$ws.Range("L40").Value = 2746.1667
$ws.Range("M40").Value = -1765.091
$ws.Range("N40").Value = -3018.1667
$ws.Range("H46").Value = 1473.625
$ws.Range("I46").Value = 1719.8
$ws.Range("J46").Value = 1063.3334
$ws.Range("K46").Value = 1719.8
$ws.Range("L46").Value = 1063.3334
$ws.Range("M46").Value = -1531.8
$ws.Range("N46").Value = -1439.3334
$ws.Range("H68").Value = 15934
$ws.Range("I68").Value = 14998
$ws.Range("J68").Value = 16150
$ws.Range("K68").Value = 14998
$ws.Range("L68").Value = 16150
$ws.Range("M68").Value = -14249
$ws.Range("N68").Value = -17648
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H71").Value = 15934
$ws.Range("I71").Value = 14998
$ws.Range("J71").Value = 16150
$ws.Range("K71").Value = 74990
$ws.Range("L71").Value = 80750
$ws.Range("M71").Value = -71246
$ws.Range("N71").Value = -88238
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H93").Value = 2933.6667
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H122").Value = 5110.3
$ws.Range("I122").Value = 6665
$ws.Range("K122").Value = 19995
$ws.Range("M122").Value = -17545
$ws.Range("H126").Value = 4051.5715
$ws.Range("I126").Value = 4051.5715
$ws.Range("K126").Value = 12154.7145
$ws.Range("M126").Value = -9684.7145
$ws.Range("H131").Value = 68996.5
$ws.Range("J131").Value = 68996.5
$ws.Range("L131").Value = 68996.5
$ws.Range("N131").Value = -79076.5
$ws.Range("H132").Value = 17863606
$ws.Range("I132").Value = 20839208
$ws.Range("K132").Value = 62517624
$ws.Range("M132").Value = -62515094

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3565.9524
$ws.Range("I96").Value = 1399.5
$ws.Range("K96").Value = 1399.5
$ws.Range("M96").Value = -26.5
$ws.Range("H122").Value = 3042.1428
$ws.Range("I122").Value = 2466.3333
$ws.Range("J122").Value = 3474
$ws.Range("K122").Value = 7398.999899999999
$ws.Range("L122").Value = 10422
$ws.Range("M122").Value = -4948.999899999999
$ws.Range("N122").Value = -15322
$ws.Range("H132").Value = 16670037
